$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INA + VGA")
$wsBattery = $wb.Worksheets.Item("M-Chain Battery Estimate")

# --- Read current (pre-sort) values of the two input columns (A,B) for rows 2..7 ---
$a2 = $ws.Cells.Item(2,1).Value()
$b2 = $ws.Cells.Item(2,2).Value()
$a3 = $ws.Cells.Item(3,1).Value()
$b3 = $ws.Cells.Item(3,2).Value()
$a4 = $ws.Cells.Item(4,1).Value()
$b4 = $ws.Cells.Item(4,2).Value()
$a5 = $ws.Cells.Item(5,1).Value()
$b5 = $ws.Cells.Item(5,2).Value()
$a6 = $ws.Cells.Item(6,1).Value()
$b6 = $ws.Cells.Item(6,2).Value()
$a7 = $ws.Cells.Item(7,1).Value()
$b7 = $ws.Cells.Item(7,2).Value()

# --- Write back the values sorted ascending by column A (mirrors a Data > Sort A-Z on A2:K7) ---
# new row2 <- old row4 (0.15)
$ws.Cells.Item(2,1).Value = $a4
$ws.Cells.Item(2,2).Value = $b4
# new row3 <- old row6 (2)
$ws.Cells.Item(3,1).Value = $a6
$ws.Cells.Item(3,2).Value = $b6
# new row4 <- old row5 (20)
$ws.Cells.Item(4,1).Value = $a5
$ws.Cells.Item(4,2).Value = $b5
# new row5 <- old row7 (45)
$ws.Cells.Item(5,1).Value = $a7
$ws.Cells.Item(5,2).Value = $b7
# new row6 <- old row3 (250)
$ws.Cells.Item(6,1).Value = $a3
$ws.Cells.Item(6,2).Value = $b3
# new row7 <- old row2 (3000)
$ws.Cells.Item(7,1).Value = $a2
$ws.Cells.Item(7,2).Value = $b2

# Record the sort state so the worksheet reflects a genuine Sort operation on A2:K7 by column A
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A7"))
$sortObj.SetRange($ws.Range("A2:K7"))
$sortObj.Header = 0
$sortObj.Apply()

# --- Add the missing "average of Saturated Output" cell in C8, matching the existing D8 average ---
$ws.Range("C8").Formula = "=AVERAGE(C2:C7)"
$ws.Range("C8").WrapText = $true
$ws.Range("C8").VerticalAlignment = -4108

# --- Apply 2-decimal ("0.00") number formatting, matching the commit "Sig fig in Excel" ---
# Order matters for style index allocation: F:G first (index 6), then D (index 7)
$ws.Range("F2:G7").NumberFormat = "0.00"
$ws.Range("D2:D7").NumberFormat = "0.00"

# --- Update the active sheet / selection to match the edited workbook state ---
[void]$wsBattery.Range("H6").Select()
[void]$ws.Activate()
[void]$ws.Range("H2").Select()
